$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: a PTY005_PartyDetailsEnquiryReject test case, cloned from row 4
# (PTY004_PartyDetailsEnquirySearch) with updated rowid/Test_Case keywords.

# New keyword strings used by the onboarding-reject test case
$ws.Range("A5").Value = "'4"
$ws.Range("B5").Value = "PTY005_PartyDetailsEnquiryReject"

# Remaining values mirror row 4 (PTY004_PartyDetailsEnquirySearch) exactly
$ws.Range("C5").Value  = "LocalTechnology Ltd 4342124"
$ws.Range("D5").Value  = 1414849
$ws.Range("E5").Value  = "LT Ltd"
$ws.Range("F5").Value  = "Local Private"
$ws.Range("G5").Value  = "'00000001"
$ws.Range("H5").Value  = "Enterprise"
$ws.Range("I5").Value  = "Australian Government, Fed Govt Dept or Auth"
$ws.Range("J5").Value  = "Full"
$ws.Range("K5").Value  = "'9890001"
$ws.Range("L5").Value  = "Australia"
$ws.Range("M5").Value  = "'2015-04-06"
$ws.Range("N5").Value  = "Australia"
$ws.Range("O5").Value  = "Australia"
$ws.Range("P5").Value  = "Electricity, Gas and Water Supply"
$ws.Range("Q5").Value  = "Gas Supply"
$ws.Range("R5").Value  = "'True"
$ws.Range("S5").Value  = "'True"
$ws.Range("T5").Value  = "'20000200001"
$ws.Range("U5").Value  = "Legal Address"
$ws.Range("V5").Value  = "Australia"
$ws.Range("W5").Value  = 2020
$ws.Range("X5").Value  = "22 SYDNEY ST"
$ws.Range("Y5").Value  = "20 SHIRLEY ST"
$ws.Range("Z5").Value  = "Pimpama Qld 30"
$ws.Range("AA5").Value = "Pennyroyal 40"
$ws.Range("AB5").Value = "Melbuorne"
$ws.Range("AC5").Value = "Victoria"
$ws.Range("AD5").Value = "Documents not collected"
$ws.Range("AE5").Value = "Party Details Enquiry"
$ws.Range("AF5").Value = "Australia"
$ws.Range("AG5").Value = "Commonwealth Bank of Australia - DBU"
$ws.Range("AH5").Value = "AU-AU"
$ws.Range("AI5").Value = "Commonwealth Bank of Australia - DBU"
$ws.Range("AJ5").Value = "Inner Works Inc 1414849"
$ws.Range("AK5").Value = "LT Ltd 4342124"
$ws.Range("AM5").Value = "LocalTechnology Ltd 4342124"
$ws.Range("AO5").Value = 9890001
$ws.Range("AP5").Value = 67855555245
$ws.Range("AQ5").Value = "Commercial Lending"
$ws.Range("AR5").Value = 1
$ws.Range("AS5").Value = 0
$ws.Range("AT5").Value = 0
$ws.Range("AU5").Value = "Subsidiary/Branch"
$ws.Range("AV5").Value = 0
$ws.Range("AW5").Value = "amipac"

# Re-apply row 4's formatting onto row 5 so the new row matches the sheet's
# established look (this also clears the quote-prefix styling that Excel
# applies to the text-forced numeric-looking values set above).
$ws.Range("A4:AW4").Copy() | Out-Null
$ws.Range("A5:AW5").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Match row 4's explicit row height
$ws.Rows.Item(5).RowHeight = 13.5

# Update the active selection to match the edited cell
$ws.Range("B5").Select() | Out-Null

Write-Host "done"
